{"js": "// Remove the trailing \"empty paragraph / page-break paragraph / copyright\n// paragraph\" trio that used to follow the last requirement line\n// (\"LOQ4095: ...\"), leaving the final empty + page-break paragraphs (which\n// stay right before the section break) untouched.\nconst searchResults = context.document.body.search(\"LOQ4095\", { matchCase: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the \"LOQ4095\" anchor paragraph.');\n}\n\n// Paragraph that contains the \"LOQ4095: ...\" requirement line.\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\n\n// The three paragraphs right after it are the ones being dropped:\n//   1) an empty \"Normal\" paragraph\n//   2) an empty \"Normal\" paragraph with pageBreakBefore\n//   3) the \"\u00a9 2020 ...\" copyright paragraph\nconst p1 = anchorParagraph.getNext();\nconst p2 = p1.getNext();\nconst p3 = p2.getNext();\n\np3.load(\"text\");\nawait context.sync();\n\n// Sanity-check we are about to remove the right paragraph before mutating.\nif (p3.text.indexOf(\"2020\") === -1) {\n  throw new Error(\"Unexpected document structure \u2014 aborting deletion.\");\n}\n\n// Delete back-to-front so earlier references stay valid.\np3.delete();\np2.delete();\np1.delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"empty paragraph / page-break paragraph / copyright\n# paragraph\" trio that used to follow the last requirement line\n# (\"LOQ4095: ...\"), leaving the final empty + page-break paragraphs (which\n# stay right before the section break) untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph containing the \"LOQ4095\" requirement line.\n$anchor = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $candidate = $d.Paragraphs($i)\n    if ($candidate.Range.Text -like \"*LOQ4095*\") {\n        $anchor = $candidate\n        break\n    }\n}\n\nif ($null -eq $anchor) {\n    throw 'Could not find the \"LOQ4095\" anchor paragraph.'\n}\n\n# The three paragraphs right after it are the ones being dropped:\n#   1) an empty \"Normal\" paragraph\n#   2) an empty \"Normal\" paragraph with pageBreakBefore\n#   3) the \"(c) 2020 ...\" copyright paragraph\n$p1 = $anchor.Next()\n$p2 = $p1.Next()\n$p3 = $p2.Next()\n\nif ($p3.Range.Text -notlike \"*2020*\") {\n    throw \"Unexpected document structure - aborting deletion.\"\n}\n\n$start = $p1.Range.Start\n$end = $p3.Range.End\n$deleteRange = $d.Range($start, $end)\n$deleteRange.Delete()\n"}
